# Auto-generated edit script applying the diff to 杭州-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 72
$ws.Range('F3').Value = 43
$ws.Range('F4').Value = 447
$ws.Range('F5').Value = 34
$ws.Range('F6').Value = 337
$ws.Range('F7').Value = 381
$ws.Range('F8').Value = 71
$ws.Range('F9').Value = 69
$ws.Range('F10').Value = 34
$ws.Range('F11').Value = 685
$ws.Range('F12').Value = 1520
$ws.Range('F13').Value = 5877
$ws.Range('F15').Value = 1655
$ws.Range('F16').Value = 406
$ws.Range('F17').Value = 5606
$ws.Range('F18').Value = 102
$ws.Range('F19').Value = 45
$ws.Range('F20').Value = 142
$ws.Range('F21').Value = 91
$ws.Range('F22').Value = 1584
$ws.Range('F23').Value = 830
$ws.Range('F24').Value = 33
$ws.Range('F25').Value = 101
$ws.Range('F26').Value = 1191
$ws.Range('G26').Value = 69
$ws.Range('F27').Value = 686
$ws.Range('F28').Value = 163
$ws.Range('F30').Value = 72
$ws.Range('F31').Value = 3839

$ws = $wb.Worksheets.Item('演出')
$ws.Range('F3').Value = 19
$ws.Range('F4').Value = 107
$ws.Range('F5').Value = 211
$ws.Range('F8').Value = 320

$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F2').Value = 9463
$ws.Range('F4').Value = 2186
$ws.Range('F5').Value = 553

$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 9463
$ws.Range('F4').Value = 2186
$ws.Range('F5').Value = 72
$ws.Range('F6').Value = 447
$ws.Range('F7').Value = 553
$ws.Range('F8').Value = 337
$ws.Range('F9').Value = 381
$ws.Range('F10').Value = 71
$ws.Range('F11').Value = 69
$ws.Range('C13').Value = '杭州·申放送x三月兽mini联动   特别健康cafe'
$ws.Range('D13').Value = '杭州in77店D区B2层B2007室 三月兽mini杭州店'
$ws.Range('E13').Value = '2024.09.07 11:00-09.08 18:10'
$ws.Range('F13').Value = 34
$ws.Range('G13').Value = 199
$ws.Range('H13').Value = 'https://show.bilibili.com/platform/detail.html?id=91134'
$ws.Range('I13').Value = '//i0.hdslb.com/bfs/openplatform/202408/L6PmCJhd1724324086367.png'
$ws.Range('B14').Value = '''2024-09-07'
$ws.Range('C14').Value = '杭州·红楼梦·主题演绎国风音乐会《梦寻红楼》'
$ws.Range('D14').Value = '望梅路与汀兰路交叉口向南100米 杭州临平大剧院（原余杭大剧院）'
$ws.Range('E14').Value = '2024.09.07 15:00-09.07 16:30'
$ws.Range('F14').Value = 19
$ws.Range('G14').Value = 100
$ws.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=89257'
$ws.Range('I14').Value = '//i2.hdslb.com/bfs/openplatform/202407/tkm6AHo71720572975141.jpeg'
$ws.Range('C15').Value = '杭州·2024首届COMIC GALAXY次元盛典'
$ws.Range('D15').Value = '长江南路336号 白马湖国际会展中心'
$ws.Range('E15').Value = '2024.09.15 09:30-09.17 17:30'
$ws.Range('F15').Value = 685
$ws.Range('G15').Value = 88
$ws.Range('H15').Value = 'https://show.bilibili.com/platform/detail.html?id=90433'
$ws.Range('I15').Value = '//i0.hdslb.com/bfs/openplatform/202408/teoBMbzd1723019674766.png'
$ws.Range('C16').Value = '杭州·浮游猫动漫嘉年华'
$ws.Range('D16').Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws.Range('E16').Value = '2024.09.15 09:00-09.16 18:00'
$ws.Range('F16').Value = 1520
$ws.Range('H16').Value = 'https://show.bilibili.com/platform/detail.html?id=88498'
$ws.Range('I16').Value = '//i2.hdslb.com/bfs/openplatform/202406/qsuFy4iv1719569431608.jpeg'
$ws.Range('C17').Value = '杭州·理想乡动漫展-同人创作者大会'
$ws.Range('D17').Value = '阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心'
$ws.Range('E17').Value = '2024.09.15 10:00-09.16 17:00'
$ws.Range('F17').Value = 5877
$ws.Range('G17').Value = 68
$ws.Range('H17').Value = 'https://show.bilibili.com/platform/detail.html?id=83822'
$ws.Range('I17').Value = '//i2.hdslb.com/bfs/openplatform/202408/oJHXYwDF1722916275016.jpeg'
$ws.Range('B18').Value = '''2024-09-15'
$ws.Range('C18').Value = '杭州·西溪银泰 布谷布Goods二次元吃谷嘉年华 免票'
$ws.Range('D18').Value = '双龙街588号 西溪银泰城'
$ws.Range('E18').Value = '2024.09.15 10:00-09.17 20:00'
$ws.Range('F18').Value = 96
$ws.Range('G18').Value = 30
$ws.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=89782'
$ws.Range('I18').Value = '//i0.hdslb.com/bfs/openplatform/202407/iWnJHkey1721737159663.png'
$ws.Range('F19').Value = 1655
$ws.Range('F22').Value = 406
$ws.Range('F25').Value = 5606
$ws.Range('F26').Value = 103
$ws.Range('F27').Value = 45
$ws.Range('F28').Value = 142
$ws.Range('F29').Value = 91
$ws.Range('F30').Value = 1584
$ws.Range('F31').Value = 830
$ws.Range('F32').Value = 33
$ws.Range('F33').Value = 101
$ws.Range('F34').Value = 1191
$ws.Range('G34').Value = 69
$ws.Range('F35').Value = 686
$ws.Range('F36').Value = 164
$ws.Range('F43').Value = 72
$ws.Range('F45').Value = 3839
